$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C (doble) and D (suma acumulada) for rows 3..12
foreach ($r in 3..12) {
    $ws.Range("C${r}").Formula = "=SUM(B${r},B${r})"
    $ws.Range("D${r}").Formula = "=SUM(B${r}:C${r})"
}

# Center-align the first new formula cell (C3), matching the style added in cellXfs
$ws.Range("C3").HorizontalAlignment = -4108   # xlCenter

# Row 14: averages per column
$ws.Range("B14").Formula = "=AVERAGE(B3:B12)"
$ws.Range("C14").Formula = "=AVERAGE(C3:C12)"
$ws.Range("D14").Formula = "=AVERAGE(D3:D12)"

# Row 17: totals (array formulas broadcasting the row-14 average across the column range)
$ws.Range("B17").FormulaArray = "=SUM(D14+B3:B12)"
$ws.Range("C17").FormulaArray = "=SUM(E14+C3:C12)"
$ws.Range("D17").FormulaArray = "=SUM(F14+D3:D12)"

# Update the active selection
[void]$ws.Range("E3").Select()
